$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 0.91
$ws.Cells.Item(2, 4).Value = 4.142
$ws.Cells.Item(2, 5).Value = 1.819
$ws.Cells.Item(2, 6).Value = 1.484
$ws.Cells.Item(2, 7).Value = 1.628
$ws.Cells.Item(2, 8).Value = 1.915
$ws.Cells.Item(2, 9).Value = 1.077
$ws.Cells.Item(2, 10).Value = 2.466
$ws.Cells.Item(2, 11).Value = 2.633
$ws.Cells.Item(2, 12).Value = 1.939
$ws.Cells.Item(2, 13).Value = 1.891
$ws.Cells.Item(2, 14).Value = 2.203
$ws.Cells.Item(2, 15).Value = 24.107
$ws.Cells.Item(3, 3).Value = 1.436
$ws.Cells.Item(3, 4).Value = 9.911
$ws.Cells.Item(3, 5).Value = 3.998
$ws.Cells.Item(3, 6).Value = 4.908
$ws.Cells.Item(3, 7).Value = 7.565
$ws.Cells.Item(3, 8).Value = 5.195
$ws.Cells.Item(3, 9).Value = 1.963
$ws.Cells.Item(3, 10).Value = 4.86
$ws.Cells.Item(3, 11).Value = 10.725
$ws.Cells.Item(3, 12).Value = 11.108
$ws.Cells.Item(3, 13).Value = 4.668
$ws.Cells.Item(3, 14).Value = 4.309
$ws.Cells.Item(3, 15).Value = 70.646
$ws.Cells.Item(4, 3).Value = 0.024
$ws.Cells.Item(4, 4).Value = 0.67
$ws.Cells.Item(4, 5).Value = 0.12
$ws.Cells.Item(4, 6).Value = 0.024
$ws.Cells.Item(4, 7).Value = 0.335
$ws.Cells.Item(4, 8).Value = 0.335
$ws.Cells.Item(4, 9).Value = 0.12
$ws.Cells.Item(4, 10).Value = 0.287
$ws.Cells.Item(4, 11).Value = 0.335
$ws.Cells.Item(4, 12).Value = 0.096
$ws.Cells.Item(4, 13).Value = 0.12
$ws.Cells.Item(4, 14).Value = 0.168
$ws.Cells.Item(4, 15).Value = 2.634
$ws.Cells.Item(5, 3).Value = 0.192
$ws.Cells.Item(5, 4).Value = 0.168
$ws.Cells.Item(5, 5).Value = 0.263
$ws.Cells.Item(5, 6).Value = 0.192
$ws.Cells.Item(5, 7).Value = 0.383
$ws.Cells.Item(5, 8).Value = 0.144
$ws.Cells.Item(5, 9).Value = 0.096
$ws.Cells.Item(5, 10).Value = 0.168
$ws.Cells.Item(5, 11).Value = 0.407
$ws.Cells.Item(5, 12).Value = 0.311
$ws.Cells.Item(5, 13).Value = 0.096
$ws.Cells.Item(5, 14).Value = 0.192
$ws.Cells.Item(5, 15).Value = 2.612
$ws.Cells.Item(6, 3).Value = 2.562
$ws.Cells.Item(6, 4).Value = 14.891
$ws.Cells.Item(6, 5).Value = 6.2
$ws.Cells.Item(6, 6).Value = 6.608000000000001
$ws.Cells.Item(6, 7).Value = 9.911
$ws.Cells.Item(6, 8).Value = 7.589
$ws.Cells.Item(6, 9).Value = 3.256
$ws.Cells.Item(6, 10).Value = 7.781000000000001
$ws.Cells.Item(6, 11).Value = 14.1
$ws.Cells.Item(6, 12).Value = 13.454
$ws.Cells.Item(6, 13).Value = 6.775
$ws.Cells.Item(6, 14).Value = 6.872
$ws.Cells.Item(6, 15).Value = 99.999
